# Auto-generated Excel COM-interop script to apply the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.352.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.095.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.49%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.69%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.089.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.81%  "

$ws.Range("E10").Value = "  -5.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.43"
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = "  -7.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.595.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.348.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.61%  "

$ws.Range("E17").Value = "  -3.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.092.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "492.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.05%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.12%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.73%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.52%  "

$ws.Range("E32").Value = "  -3.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "520.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.71%  "

$ws.Range("E36").Value = "  -6.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0408"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.141.77"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0806"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.120"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.258"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.29%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.109"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.81%  "

